$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the results columns (C:H) for the data rows 2-15; these previously
# held per-row lookup results (stato, protocollo uscita, provvedimento,
# data provvedimento, codice richiesta, note usmaf) that are no longer
# populated. Column A (AWB) and column B (ricerca timestamp) are left as-is.
$ws.Range("C2:H15").Value = $null

# Update the active selection left by the editing session.
[void]$ws.Range("F18").Select()

$wb.Save()
